# fix French translation for iso27005 matrix
#
# - bump library_version (library_content!B2): 2 -> 3
# - fix French "likelihood" labels on the spec sheet:
#     L5: "2 - plutôt probable"  -> "2 - plutôt improbable"
#     L6: "1 - peu probable"     -> "1 - improbable"
# - widen column L on the spec sheet so the new, longer French label fits
# - restore selection/active-sheet state left over from editing

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("library_content")
$ws2 = $wb.Worksheets.Item("spec")

# Bump the library version to reflect the content fix.
$ws1.Range("B2").Value = 3

# Correct the French "likelihood" wording (these were swapped/mistranslated).
$ws2.Range("L5").Value = "2 - plutôt improbable"
$ws2.Range("L6").Value = "1 - improbable"

# The new French text is longer, so widen column L to fit it.
$ws2.Columns.Item(12).ColumnWidth = 17.2

# Leave the selection on library_content (first sheet), matching where the
# author ended up after making the fix.
$ws1.Activate()
$ws1.Range("B3").Select()

$ws2.Activate()
$ws2.Range("L6").Select()
$ws1.Activate()
